$wb = $excel.ActiveWorkbook

# The "changelog" sheet is the 3rd sheet in the workbook.
$changelog = $wb.Worksheets.Item("changelog")

# Add a new changelog row (row 5): ID = 4, Date = 2025-08-12 (serial 45881)
$changelog.Range("A5").Value = 4
$changelog.Range("B5").Value = (Get-Date -Year 2025 -Month 8 -Day 12 -Hour 0 -Minute 0 -Second 0)

# Update the selection on the changelog sheet to C5
$changelog.Range("C5").Select() | Out-Null

# Activate the changelog sheet, making it the active/selected tab
$changelog.Activate()

$wb.Save()
